# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 108600571
$ws.Range("B3").Value = 8377
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 106545
$ws.Range("F3").Value = 'Mindre märgborre'
$ws.Range("G3").Value = 'Tomicus minor'
$ws.Range("H3").Value = '(Hartig, 1834)'
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("M3").Value = 'äldre gnagspår'
$ws.Range("Q3").Value = 528457.9539003669
$ws.Range("R3").Value = 6541148.642409162

# Row 4
$ws.Range("A4").Value = 108600717
$ws.Range("B4").Value = 89412
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 5442
$ws.Range("F4").Value = 'Tallticka'
$ws.Range("G4").Value = 'Porodaedalea pini'
$ws.Range("H4").Value = '(Brot.) Murrill'
$ws.Range("I4").Value = "'2"
$ws.Range("J4").Value = 'fruktkroppar'
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 528514.2341294074
$ws.Range("R4").Value = 6541351.991765272

# Row 5
$ws.Range("A5").Value = 108600671
$ws.Range("B5").Value = 8377
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 106545
$ws.Range("F5").Value = 'Mindre märgborre'
$ws.Range("G5").Value = 'Tomicus minor'
$ws.Range("H5").Value = '(Hartig, 1834)'
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("M5").Value = 'äldre gnagspår'
$ws.Range("Q5").Value = 528575.3995984152
$ws.Range("R5").Value = 6541244.801401596

# Row 6
$ws.Range("A6").Value = 108600691
$ws.Range("B6").Value = 78098
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 6453
$ws.Range("F6").Value = 'Vedskivlav'
$ws.Range("G6").Value = 'Hertelidea botryosa'
$ws.Range("H6").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = 'bålar'
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value = 528584.0134725916
$ws.Range("R6").Value = 6541265.982774138

# Row 7
$ws.Range("A7").Value = 108600554
$ws.Range("B7").Value = 8377
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 106545
$ws.Range("F7").Value = 'Mindre märgborre'
$ws.Range("G7").Value = 'Tomicus minor'
$ws.Range("H7").Value = '(Hartig, 1834)'
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("M7").Value = 'äldre gnagspår'
$ws.Range("Q7").Value = 528511.5568710293
$ws.Range("R7").Value = 6541087.231412024

# Row 8
$ws.Range("A8").Value = 108600864
$ws.Range("B8").Value = 5426
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 101410
$ws.Range("F8").Value = 'Reliktbock'
$ws.Range("G8").Value = 'Nothorhina muricata'
$ws.Range("H8").Value = '(Dalman, 1817)'
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("M8").Value = 'färska gnagspår'
$ws.Range("Q8").Value = 528502.4862113618
$ws.Range("R8").Value = 6541335.938193527

# Row 9
$ws.Range("A9").Value = 108600622
$ws.Range("B9").Value = 89412
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 5442
$ws.Range("F9").Value = 'Tallticka'
$ws.Range("G9").Value = 'Porodaedalea pini'
$ws.Range("H9").Value = '(Brot.) Murrill'
$ws.Range("I9").Value = "'1"
$ws.Range("J9").Value = 'fruktkroppar'
$ws.Range("M9").ClearContents()
$ws.Range("Q9").Value = 528458.0765920902
$ws.Range("R9").Value = 6541201.695022714

# Row 10
$ws.Range("A10").Value = 108600779
$ws.Range("B10").Value = 8377
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 106545
$ws.Range("F10").Value = 'Mindre märgborre'
$ws.Range("G10").Value = 'Tomicus minor'
$ws.Range("H10").Value = '(Hartig, 1834)'
$ws.Range("I10").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("M10").Value = 'äldre gnagspår'
$ws.Range("Q10").Value = 528499.9212518559
$ws.Range("R10").Value = 6541472.918336567

# Row 11
$ws.Range("A11").Value = 108600651
$ws.Range("B11").Value = 8367
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 106554
$ws.Range("F11").Value = 'Björksplintborre'
$ws.Range("G11").Value = 'Scolytus ratzeburgii'
$ws.Range("H11").Value = 'Janson, 1856'
$ws.Range("I11").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("M11").Value = 'färska gnagspår'
$ws.Range("Q11").Value = 528559.4282089664
$ws.Range("R11").Value = 6541241.59220575

# Row 12
$ws.Range("A12").Value = 108600738
$ws.Range("B12").Value = 78098
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 6453
$ws.Range("F12").Value = 'Vedskivlav'
$ws.Range("G12").Value = 'Hertelidea botryosa'
$ws.Range("H12").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("I12").ClearContents()
$ws.Range("J12").Value = 'bålar'
$ws.Range("M12").ClearContents()
$ws.Range("Q12").Value = 528501.4884462073
$ws.Range("R12").Value = 6541400.825744567

# Row 13
$ws.Range("A13").Value = 108600599
$ws.Range("B13").Value = 8377
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 106545
$ws.Range("F13").Value = 'Mindre märgborre'
$ws.Range("G13").Value = 'Tomicus minor'
$ws.Range("H13").Value = '(Hartig, 1834)'
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("M13").Value = 'äldre gnagspår'
$ws.Range("Q13").Value = 528440.9271998855
$ws.Range("R13").Value = 6541148.516224748

# Row 14
$ws.Range("A14").Value = 108600803
$ws.Range("B14").Value = 8367
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 106554
$ws.Range("F14").Value = 'Björksplintborre'
$ws.Range("G14").Value = 'Scolytus ratzeburgii'
$ws.Range("H14").Value = 'Janson, 1856'
$ws.Range("I14").ClearContents()
$ws.Range("J14").ClearContents()
$ws.Range("M14").Value = 'färska gnagspår'
$ws.Range("Q14").Value = 528455.3571655933
$ws.Range("R14").Value = 6541429.325337943
